# Update the dSF column (column F) values to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = -6
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = -1
